# Fruta / hortaliza, semanal
# Insert 2 new weekly records (rows) for "Vega Monumental Concepción - Ciruela"
# right before the existing row 138, pushing the old rows 138-142 down to 140-144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 138 (old rows 138:139 shift down to 140:141, etc.)
$ws.Range("A138:T139").EntireRow.Insert()

# --- New row 138 ---
$ws.Cells.Item(138, 1).Value = 11
$ws.Cells.Item(138, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(138, 3).Value = "Bíobío"
$ws.Cells.Item(138, 4).Value = 45041
$ws.Cells.Item(138, 5).Value = 8
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100103
$ws.Cells.Item(138, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(138, 9).Value = 100103002
$ws.Cells.Item(138, 10).Value = "Ciruela"
$ws.Cells.Item(138, 11).Value = "Angeleno"
$ws.Cells.Item(138, 12).Value = "Primera"
$ws.Cells.Item(138, 13).Value = 50
$ws.Cells.Item(138, 14).Value = 10000
$ws.Cells.Item(138, 15).Value = 10000
$ws.Cells.Item(138, 16).Value = 10000
$ws.Cells.Item(138, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(138, 18).Value = "Región del Maule"
$ws.Cells.Item(138, 19).Value = 556
$ws.Cells.Item(138, 20).Value = 18

# --- New row 139 ---
$ws.Cells.Item(139, 1).Value = 11
$ws.Cells.Item(139, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(139, 3).Value = "Bíobío"
$ws.Cells.Item(139, 4).Value = 45041
$ws.Cells.Item(139, 5).Value = 8
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100103
$ws.Cells.Item(139, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(139, 9).Value = 100103002
$ws.Cells.Item(139, 10).Value = "Ciruela"
$ws.Cells.Item(139, 11).Value = "Angeleno"
$ws.Cells.Item(139, 12).Value = "Segunda"
$ws.Cells.Item(139, 13).Value = 50
$ws.Cells.Item(139, 14).Value = 9000
$ws.Cells.Item(139, 15).Value = 9000
$ws.Cells.Item(139, 16).Value = 9000
$ws.Cells.Item(139, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(139, 18).Value = "Región del Maule"
$ws.Cells.Item(139, 19).Value = 500
$ws.Cells.Item(139, 20).Value = 18

# Make sure the date cells keep the date number format used elsewhere in column D
$ws.Range("D138:D139").NumberFormat = $ws.Range("D137").NumberFormat

Write-Host "Dimension now:" $ws.UsedRange.Address()
